$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the existing row 396 (Andross/Carson block),
# pushing the current rows 396-405 down to 400-409.
$ws.Rows("396:399").Insert()

# Common columns shared by every data row in this block.
$commonA = 5
$commonB = "Macroferia Regional de Talca"
$commonC = "Maule"
$commonE = 7
$commonF = "Fruta"
$commonG = 100103
$commonH = "Frutos de hueso (carozo)"
$commonI = 100103004
$commonJ = "Durazno"

function Set-Row {
    param($r, $d, $k, $l, $m, $n, $o, $p, $q, $rOrigen, $s, $t)

    $ws.Cells.Item($r, 1).Value = $commonA
    $ws.Cells.Item($r, 2).Value = $commonB
    $ws.Cells.Item($r, 3).Value = $commonC
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $commonE
    $ws.Cells.Item($r, 6).Value = $commonF
    $ws.Cells.Item($r, 7).Value = $commonG
    $ws.Cells.Item($r, 8).Value = $commonH
    $ws.Cells.Item($r, 9).Value = $commonI
    $ws.Cells.Item($r, 10).Value = $commonJ
    $ws.Cells.Item($r, 11).Value = $k
    $ws.Cells.Item($r, 12).Value = $l
    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 14).Value = $n
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p
    $ws.Cells.Item($r, 17).Value = $q
    $ws.Cells.Item($r, 18).Value = $rOrigen
    $ws.Cells.Item($r, 19).Value = $s
    $ws.Cells.Item($r, 20).Value = $t
}

# New row 396: Early Majestic / Primera
Set-Row 396 44890 "Early Majestic" "Primera" 230 18000 18000 18000 "$/caja 16 kilos granel" "Región de O'Higgins" 1125 16

# New row 397: Early Majestic / Segunda
Set-Row 397 44890 "Early Majestic" "Segunda" 180 16000 16000 16000 "$/caja 16 kilos granel" "Región de O'Higgins" 1000 16

# New row 398: Florida King / Primera (note origin region changes to Coquimbo)
Set-Row 398 44890 "Florida King" "Primera" 120 16000 16000 16000 "$/caja 16 kilos granel" "Región de Coquimbo" 1000 16

# New row 399: Florida King / Segunda
Set-Row 399 44890 "Florida King" "Segunda" 100 14000 14000 14000 "$/caja 16 kilos granel" "Región de O'Higgins" 875 16
